$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (column C) date value from 45202 to 45203
# for rows 2 through 18.
for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
